# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scheduled
# update). Most cells are plain text (Price/Volume columns are stored as
# literal strings, not numbers, e.g. "26.271.18" / "  +0.15%  "), and two
# rows (39/40) swap which coin ("Maker" <-> "VeChain") they describe.
#
# Price strings that parse as a plain decimal (e.g. "1.008") would
# otherwise be auto-coerced to a Number by Excel's normal type inference
# when assigned through Range.Value, silently dropping the original text
# formatting (trailing zeros, etc.) and flipping the cell's stored type.
# Prefixing the literal with an apostrophe is exactly what typing such a
# value into Excel and forcing "Text" entry does: Excel stores the text
# verbatim (dropping only the leading apostrophe marker) and keeps the
# General number format, matching the source values exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.271.18'
$ws.Range('E2').Value = '  +0.15%  '

$ws.Range('D3').Value = '1.679.13'
$ws.Range('E3').Value = '  +0.39%  '

$ws.Range('D4').Value = '''1.008'
$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = '''218.29'
$ws.Range('E5').Value = '  +0.23%  '

$ws.Range('D6').Value = '''0.5258'
$ws.Range('E6').Value = '  +2.66%  '

$ws.Range('E7').Value = '  +0.09%  '

$ws.Range('D8').Value = '''0.2697'
$ws.Range('E8').Value = '  +1.32%  '

$ws.Range('D9').Value = '''0.06432'
$ws.Range('E9').Value = '  +0.77%  '

$ws.Range('E10').Value = '  +2.00%  '

$ws.Range('D11').Value = '''0.07502'
$ws.Range('E11').Value = '  +1.52%  '

$ws.Range('D12').Value = '1.694.67'
$ws.Range('E12').Value = '  +1.29%  '

$ws.Range('E13').Value = '  -0.29%  '

$ws.Range('D14').Value = '''0.5806'
$ws.Range('E14').Value = '  -0.49%  '

$ws.Range('D15').Value = '''0.000008480'
$ws.Range('E15').Value = '  -2.07%  '

$ws.Range('D16').Value = '''64.25'
$ws.Range('E16').Value = '  -0.56%  '

$ws.Range('D17').Value = '26.309.74'
$ws.Range('E17').Value = '  -0.03%  '

$ws.Range('D18').Value = '''4.925'
$ws.Range('E18').Value = '  -0.80%  '

$ws.Range('D19').Value = '''1.008'
$ws.Range('E19').Value = '  +0.15%  '

$ws.Range('E20').Value = '  -0.36%  '

$ws.Range('D21').Value = '''189.41'
$ws.Range('E21').Value = '  -0.08%  '

$ws.Range('D22').Value = '''6.199'
$ws.Range('E22').Value = '  -0.25%  '

$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('D24').Value = '''144.96'
$ws.Range('E24').Value = '  +0.60%  '

$ws.Range('D25').Value = '''7.724'
$ws.Range('E25').Value = '  +0.83%  '

$ws.Range('D26').Value = '''0.1240'
$ws.Range('E26').Value = '  +4.85%  '

$ws.Range('D27').Value = '''15.80'
$ws.Range('E27').Value = '  +0.78%  '

$ws.Range('D28').Value = '''0.06598'
$ws.Range('E28').Value = '  +9.63%  '

$ws.Range('D29').Value = '''1.357'
$ws.Range('E29').Value = '  +5.79%  '

$ws.Range('D30').Value = '''1.327'
$ws.Range('E30').Value = '  +0.04%  '

$ws.Range('D31').Value = '''3.583'
$ws.Range('E31').Value = '  +1.48%  '

$ws.Range('D32').Value = '''3.570'
$ws.Range('E32').Value = '  +0.93%  '

$ws.Range('D33').Value = '''1.661'
$ws.Range('E33').Value = '  +0.82%  '

$ws.Range('D34').Value = '''1.026'
$ws.Range('E34').Value = '  +0.82%  '

$ws.Range('D35').Value = '''0.6200'
$ws.Range('E35').Value = '  +2.72%  '

$ws.Range('D36').Value = '''2.397'
$ws.Range('E36').Value = '  +0.90%  '

$ws.Range('D37').Value = '''2.720'
$ws.Range('E37').Value = '  +2.76%  '

$ws.Range('D38').Value = '''6.389'

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.01621'
$ws.Range('E39').Value = '  +0.16%  '

$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '1.105.23'
$ws.Range('E40').Value = '  +2.10%  '

$ws.Range('D41').Value = '''0.8753'
$ws.Range('E41').Value = '  +0.55%  '

$ws.Range('E42').Value = '  +0.34%  '

$ws.Range('D43').Value = '''100.45'
$ws.Range('E43').Value = '  -0.05%  '

$ws.Range('D44').Value = '1.827.59'
$ws.Range('E44').Value = '  +0.34%  '

$ws.Range('D45').Value = '''0.00000000112'
$ws.Range('E45').Value = '  -0.70%  '

$ws.Range('D46').Value = '''56.82'
$ws.Range('E46').Value = '  +0.70%  '

$ws.Range('E47').Value = '  -0.48%  '

$ws.Range('D48').Value = '''8.112'
$ws.Range('E48').Value = '  +0.25%  '

$ws.Range('D49').Value = '''0.05271'
$ws.Range('E49').Value = '  +1.00%  '

$ws.Range('D50').Value = '''0.4302'

$ws.Range('D51').Value = '''6.054'
$ws.Range('E51').Value = '  +2.67%  '
